# "Generate Report for Handoff"
#
# The localization-status report is regenerated: every row whose handoff
# date previously showed an already-reported time gets refreshed to the
# new report-generation timestamps. This affects the "Latest Handoff
# Date" column on the Overview sheet (column D) and the "Latest Handoff
# Datetime" column on the per-locale sheets (column E), for the File
# Name rows: 3677d9e5...md, 32d33a1d...md, 5256e1cc...md, 669c895a...md,
# 9d786cc5...md, bd856a0b...md, d5cd070d...md, f534a21d...md
# (spreadsheet rows 7 and 10-16 on every sheet).

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

# --- Overview sheet: column D = "Latest Handoff Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 4).Value = "2016-03-22 08:27:18"
}

# --- zh-cn sheet: column E = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "2016-03-22 08:27:13"
}

# --- de-de sheet: column E = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "2016-03-22 08:27:18"
}
